# "added new teststeps from macbook"
#
# Adds a third worksheet ("AddPIMUserData") to the workbook, right after
# "AddUserData". The new sheet carries the same single header/seed row as
# "AddUserData" (Admin / admin123 / e / Jack Jonson / Test@123 / Test@123,
# including the two mailto: hyperlinks on E1/F1) - copying the sheet is the
# simplest faithful way to reproduce that row plus its hyperlinks. It then
# becomes the active sheet/tab. Selections on all three sheets are touched:
# LoginData keeps its active cell but the selected range becomes the whole
# first row; AddUserData's selection collapses back to A1 (whole row 1) and
# it is no longer the selected tab; the new sheet opens selected at A1
# (whole row 1) as the now-active tab.

$wb = $excel.ActiveWorkbook

$loginData = $wb.Worksheets.Item("LoginData")
$addUserData = $wb.Worksheets.Item("AddUserData")

# Duplicate "AddUserData" (same header row + hyperlinks) right after itself,
# then rename the copy - this becomes the new "AddPIMUserData" sheet and
# Excel makes it the active sheet/tab automatically.
$addUserData.Copy($null, $addUserData) | Out-Null
$pimUserData = $wb.Worksheets.Item($wb.Worksheets.Count)
$pimUserData.Name = "AddPIMUserData"

# LoginData: active cell stays B1 (only the selected range widens to the
# whole of row 1 in the source file).
$loginData.Activate() | Out-Null
$loginData.Range("B1").Select() | Out-Null

# AddUserData: no longer the active tab; selection resets to A1 / row 1.
$addUserData.Activate() | Out-Null
$addUserData.Rows(1).Select() | Out-Null

# AddPIMUserData: becomes (and stays) the active tab, selection at A1 / row 1.
$pimUserData.Activate() | Out-Null
$pimUserData.Rows(1).Select() | Out-Null
